# Update Sheet1's raw data block (A1:I4). Sheet2's B5:G8 cells pull these
# values via INDEX(Sheet1!...) formulas, so they recalc automatically.
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

$values = @(
    @(-53.392709111029731, -55.6188407187203, -27.9252905327618, -38.044571325913239, 49.853236446571145, 24.865290349462029, 0.71611381944429187, 0.076556153168828411, 0.017908940363758365),
    @(-35.301389668524692, -21.096290110667503, 4.3509948160231708, 1.13678158720719, 102.09694055269313, 841.61971847202858, 0.049292072964076998, 0.66716632437757439, 0.36027259808856049),
    @(-26.385158702470903, -24.939424988609087, -1.4989208996281305, 2.0893944206969048, 16.767145981410874, 85.741375495763535, 0.88959344251944805, 0.7952215365661579, 0.30395104643883963),
    @(-128.26333134283314, -96.030651582788735, -104.44950967422911, -71.222603796937065, 18.462148713722247, 28.151443458276059, 0.019075380041754363, 0.028513848259217267, 0.31846951538707147)
)

for ($r = 0; $r -lt $values.Length; $r++) {
    $row = $values[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws1.Cells.Item($r + 1, $c + 1).Value = $row[$c]
    }
}

$excel.Calculate()
